$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new cell A10 with the new comment text (will be added to shared strings)
$ws.Range("A10").Value = "#output based on sheet1 of Hacktoberfest_database.xlsx"

# Widen column A to fit new content (target stored width 52.1640625 characters;
# COM ColumnWidth rounds to whole-pixel steps, so feed the value that lands
# on the closest achievable pixel width)
$ws.Columns.Item(1).ColumnWidth = 51.333333

# Move the active selection to A12
$ws.Range("A12").Select()
